{"js": "// Update the three-digit-by-one-digit multiplication answers in the\n// practice-sheet table. Each cell holds a single \"NNN\u00d7N=NNNN\" run; we\n// locate the old equation text exactly (case/format preserving) and\n// swap in the new equation text in place, leaving all run/paragraph\n// formatting (font, size, alignment) untouched.\nconst replacements = [\n  [\"531\u00d74=2124\", \"698\u00d75=3490\"],\n  [\"586\u00d74=2344\", \"521\u00d72=1042\"],\n  [\"787\u00d73=2361\", \"337\u00d78=2696\"],\n  [\"178\u00d72=356\", \"779\u00d77=5453\"],\n  [\"947\u00d78=7576\", \"319\u00d76=1914\"],\n  [\"361\u00d79=3249\", \"855\u00d79=7695\"],\n  [\"250\u00d76=1500\", \"506\u00d76=3036\"],\n  [\"833\u00d79=7497\", \"134\u00d77=938\"],\n  [\"802\u00d76=4812\", \"120\u00d75=600\"],\n  [\"619\u00d76=3714\", \"700\u00d76=4200\"],\n  [\"484\u00d78=3872\", \"743\u00d76=4458\"],\n  [\"125\u00d79=1125\", \"734\u00d78=5872\"],\n  [\"321\u00d77=2247\", \"344\u00d73=1032\"],\n  [\"691\u00d72=1382\", \"634\u00d77=4438\"],\n  [\"543\u00d77=3801\", \"663\u00d74=2652\"],\n  [\"781\u00d77=5467\", \"148\u00d79=1332\"],\n  [\"509\u00d72=1018\", \"560\u00d78=4480\"],\n  [\"568\u00d73=1704\", \"134\u00d78=1072\"],\n  [\"667\u00d73=2001\", \"327\u00d74=1308\"],\n  [\"417\u00d78=3336\", \"542\u00d75=2710\"],\n  [\"766\u00d72=1532\", \"480\u00d74=1920\"],\n  [\"572\u00d79=5148\", \"292\u00d79=2628\"],\n  [\"402\u00d77=2814\", \"295\u00d72=590\"],\n  [\"157\u00d77=1099\", \"807\u00d79=7263\"],\n  [\"632\u00d75=3160\", \"876\u00d74=3504\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const found = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  found.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < found.items.length; i++) {\n    found.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the three-digit-by-one-digit multiplication answers in the\n# practice-sheet table. Each cell holds a single \"NNN\u00d7N=NNNN\" run; use\n# Find/Replace (whole-document range) to swap each old equation text for\n# its new value, preserving the surrounding run/paragraph formatting.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @('531\u00d74=2124', '698\u00d75=3490'),\n    @('586\u00d74=2344', '521\u00d72=1042'),\n    @('787\u00d73=2361', '337\u00d78=2696'),\n    @('178\u00d72=356',  '779\u00d77=5453'),\n    @('947\u00d78=7576', '319\u00d76=1914'),\n    @('361\u00d79=3249', '855\u00d79=7695'),\n    @('250\u00d76=1500', '506\u00d76=3036'),\n    @('833\u00d79=7497', '134\u00d77=938'),\n    @('802\u00d76=4812', '120\u00d75=600'),\n    @('619\u00d76=3714', '700\u00d76=4200'),\n    @('484\u00d78=3872', '743\u00d76=4458'),\n    @('125\u00d79=1125', '734\u00d78=5872'),\n    @('321\u00d77=2247', '344\u00d73=1032'),\n    @('691\u00d72=1382', '634\u00d77=4438'),\n    @('543\u00d77=3801', '663\u00d74=2652'),\n    @('781\u00d77=5467', '148\u00d79=1332'),\n    @('509\u00d72=1018', '560\u00d78=4480'),\n    @('568\u00d73=1704', '134\u00d78=1072'),\n    @('667\u00d73=2001', '327\u00d74=1308'),\n    @('417\u00d78=3336', '542\u00d75=2710'),\n    @('766\u00d72=1532', '480\u00d74=1920'),\n    @('572\u00d79=5148', '292\u00d79=2628'),\n    @('402\u00d77=2814', '295\u00d72=590'),\n    @('157\u00d77=1099', '807\u00d79=7263'),\n    @('632\u00d75=3160', '876\u00d74=3504')\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
